$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new test result row
$ws.Range("B5").Value = 13322880
$ws.Range("C5").Value = "ok"

# Update the active selection to D5 (matches the saved view state)
$ws.Range("D5").Select()
